$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-24 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-25 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("90×23=", $true, $false, $false, $false, $false, $true, 1, $false, "91×42=", 2) | Out-Null
$d.Content.Find.Execute("47×32=", $true, $false, $false, $false, $false, $true, 1, $false, "15×12=", 2) | Out-Null
$d.Content.Find.Execute("12×35=", $true, $false, $false, $false, $false, $true, 1, $false, "19×54=", 2) | Out-Null
$d.Content.Find.Execute("77×60=", $true, $false, $false, $false, $false, $true, 1, $false, "26×18=", 2) | Out-Null
$d.Content.Find.Execute("67×96=", $true, $false, $false, $false, $false, $true, 1, $false, "58×76=", 2) | Out-Null
$d.Content.Find.Execute("93×27=", $true, $false, $false, $false, $false, $true, 1, $false, "57×20=", 2) | Out-Null
$d.Content.Find.Execute("42×46=", $true, $false, $false, $false, $false, $true, 1, $false, "21×82=", 2) | Out-Null
$d.Content.Find.Execute("87×46=", $true, $false, $false, $false, $false, $true, 1, $false, "33×61=", 2) | Out-Null
$d.Content.Find.Execute("31×18=", $true, $false, $false, $false, $false, $true, 1, $false, "45×99=", 2) | Out-Null
$d.Content.Find.Execute("92×57=", $true, $false, $false, $false, $false, $true, 1, $false, "64×92=", 2) | Out-Null
$d.Content.Find.Execute("26×74=", $true, $false, $false, $false, $false, $true, 1, $false, "65×23=", 2) | Out-Null
$d.Content.Find.Execute("18×68=", $true, $false, $false, $false, $false, $true, 1, $false, "41×21=", 2) | Out-Null
$d.Content.Find.Execute("19×57=", $true, $false, $false, $false, $false, $true, 1, $false, "41×39=", 2) | Out-Null
$d.Content.Find.Execute("78×66=", $true, $false, $false, $false, $false, $true, 1, $false, "37×46=", 2) | Out-Null
$d.Content.Find.Execute("96×77=", $true, $false, $false, $false, $false, $true, 1, $false, "93×69=", 2) | Out-Null
$d.Content.Find.Execute("39×76=", $true, $false, $false, $false, $false, $true, 1, $false, "48×73=", 2) | Out-Null
$d.Content.Find.Execute("38×63=", $true, $false, $false, $false, $false, $true, 1, $false, "83×37=", 2) | Out-Null
$d.Content.Find.Execute("32×22=", $true, $false, $false, $false, $false, $true, 1, $false, "24×56=", 2) | Out-Null
$d.Content.Find.Execute("44×43=", $true, $false, $false, $false, $false, $true, 1, $false, "27×63=", 2) | Out-Null
$d.Content.Find.Execute("42×42=", $true, $false, $false, $false, $false, $true, 1, $false, "22×51=", 2) | Out-Null
$d.Content.Find.Execute("70×52=", $true, $false, $false, $false, $false, $true, 1, $false, "70×59=", 2) | Out-Null
$d.Content.Find.Execute("18×69=", $true, $false, $false, $false, $false, $true, 1, $false, "60×97=", 2) | Out-Null
$d.Content.Find.Execute("68×86=", $true, $false, $false, $false, $false, $true, 1, $false, "54×98=", 2) | Out-Null
$d.Content.Find.Execute("25×85=", $true, $false, $false, $false, $false, $true, 1, $false, "11×99=", 2) | Out-Null
$d.Content.Find.Execute("35×88=", $true, $false, $false, $false, $false, $true, 1, $false, "45×49=", 2) | Out-Null
